$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-empty "hours worked" cells across the three weekly
#     tables (student id column I4, plus the per-day/per-person hour grids) ---

# Row 4: student numeric id replaces the placeholder "id6" text
$ws.Range("I4").Value = 4672372

# week 3.1 totals (rows 8-11)
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 4
$ws.Range("I10").Value = 6
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 3

# week 3.2 grid (rows 14-19)
$ws.Range("D14").Value = 3

$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 4
$ws.Range("I15").Value = 4

$ws.Range("D16").Value = 4

$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 4

$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 4

$ws.Range("D19").Value = 8
$ws.Range("F19").Value = 3

# week 3.3 grid (rows 21-27)
$ws.Range("D21").Value = 4
$ws.Range("F21").Value = 4

$ws.Range("D22").Value = 8
$ws.Range("F22").Value = 6

$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 6
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 10
$ws.Range("I23").Value = 4

$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 7
$ws.Range("G24").Value = 7
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 8

$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 10
$ws.Range("I25").Value = 10

$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 12

$ws.Range("D27").Value = 8
$ws.Range("E27").Value = 6
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 3

# Selection moved to G29 in the saved file
$ws.Range("G29").Select() | Out-Null
